$p = $ppt.ActivePresentation
$s = $p.Slides.Item(27)
$s.SlideShowTransition.Hidden = $true
